$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1277.2
$ws.Range("J28").Value = 1999.5
$ws.Range("L28").Value = 1999.5
$ws.Range("N28").Value = -2969.5
$ws.Range("H33").Value = 865.25
$ws.Range("I33").Value = 44.555557
$ws.Range("K33").Value = 44.555557
$ws.Range("M33").Value = 184.444443
$ws.Range("H40").Value = 1996.6666
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 1996.6666
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 1996.6666
$ws.Range("M40").Value = ""
$ws.Range("N40").Value = -2346.6666
$ws.Range("H51").Value = 9999.286
$ws.Range("I51").Value = 10000
$ws.Range("J51").Value = 9999.166999999999
$ws.Range("K51").Value = 10000
$ws.Range("L51").Value = 9999.166999999999
$ws.Range("M51").Value = -9516
$ws.Range("N51").Value = -10967.167
$ws.Range("H70").Value = 1833.3334
$ws.Range("J70").Value = 1833.3334
$ws.Range("L70").Value = 5500.0002
$ws.Range("N70").Value = -6040.0002
$ws.Range("H73").Value = 1833.3334
$ws.Range("J73").Value = 1833.3334
$ws.Range("L73").Value = 5500.0002
$ws.Range("N73").Value = -7372.0002
$ws.Range("H132").Value = 2715.5715
$ws.Range("I132").Value = 2715.5715
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8146.7145
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5616.7145
$ws.Range("N132").Value = ""
$ws.Range("H135").Value = 2542.077
$ws.Range("I135").Value = 2599.7
$ws.Range("J135").Value = 2350
$ws.Range("K135").Value = 23397.3
$ws.Range("L135").Value = 21150
$ws.Range("M135").Value = -20862.3
$ws.Range("N135").Value = -26220
$ws.Range("H137").Value = 2465.1667
$ws.Range("J137").Value = 3000
$ws.Range("L137").Value = 9000
$ws.Range("N137").Value = -14100

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2365.4285
$ws.Range("J45").Value = 3999
$ws.Range("L45").Value = 3999
$ws.Range("N45").Value = -4753
$ws.Range("H61").Value = 3997.3333
$ws.Range("I61").Value = 3000
$ws.Range("J61").Value = 4196.8
$ws.Range("K61").Value = 3000
$ws.Range("L61").Value = 4196.8
$ws.Range("M61").Value = -2788
$ws.Range("N61").Value = -4620.8
$ws.Range("H74").Value = 1142.2858
$ws.Range("I74").Value = 1142.2858
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1142.2858
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -268.2858000000001
$ws.Range("N74").Value = ""
$ws.Range("H77").Value = 1142.2858
$ws.Range("I77").Value = 1142.2858
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 5711.429
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -1343.429
$ws.Range("N77").Value = ""
$ws.Range("H92").Value = 55000
$ws.Range("J92").Value = 55000
$ws.Range("L92").Value = 55000
$ws.Range("N92").Value = -59992
$ws.Range("H122").Value = 401
$ws.Range("I122").Value = 401
$ws.Range("K122").Value = 1203
$ws.Range("M122").Value = 1247
$ws.Range("H132").Value = 3940.3076
$ws.Range("I132").Value = 3435.3333
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 10305.9999
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -7775.999899999999
$ws.Range("N132").Value = -35060
$ws.Range("H136").Value = 3997.3333
$ws.Range("I136").Value = 3000
$ws.Range("J136").Value = 4196.8
$ws.Range("K136").Value = 9000
$ws.Range("L136").Value = 12590.4
$ws.Range("M136").Value = -6450
$ws.Range("N136").Value = -17690.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3427.6667
$ws.Range("I107").Value = 3563.2
$ws.Range("K107").Value = 3563.2
$ws.Range("M107").Value = -1643.2
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").Value = ""
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 997
$ws.Range("J22").Value = 994.5
$ws.Range("L22").Value = 994.5
$ws.Range("N22").Value = -1694.5
$ws.Range("H74").Value = 32500
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").Value = ""
$ws.Range("H77").Value = 32500
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").Value = ""
$ws.Range("H132").Value = 2794
$ws.Range("J132").Value = 4748
$ws.Range("L132").Value = 14244
$ws.Range("N132").Value = -19304

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 649.5
$ws.Range("I5").Value = 500
$ws.Range("K5").Value = 1500
$ws.Range("M5").Value = -1388
$ws.Range("H26").Value = 99
$ws.Range("J26").Value = 99
$ws.Range("L26").Value = 297
$ws.Range("N26").Value = -873
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = ""
$ws.Range("N86").Value = ""
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = ""
$ws.Range("N89").Value = ""
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").Value = ""
$ws.Range("H135").Value = 649.5
$ws.Range("I135").Value = 500
$ws.Range("K135").Value = 4500
$ws.Range("M135").Value = -1965

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 568.3125
$ws.Range("I97").Value = 549.5714
$ws.Range("J97").Value = 699.5
$ws.Range("K97").Value = 549.5714
$ws.Range("L97").Value = 699.5
$ws.Range("M97").Value = -53.57140000000004
$ws.Range("N97").Value = -1691.5
$ws.Range("H102").Value = 5555
$ws.Range("I102").Value = 5555
$ws.Range("K102").Value = 5555
$ws.Range("M102").Value = -3933
$ws.Range("H122").Value = 1952
$ws.Range("I122").Value = 2025
$ws.Range("J122").Value = 1806
$ws.Range("K122").Value = 6075
$ws.Range("L122").Value = 5418
$ws.Range("M122").Value = -3625
$ws.Range("N122").Value = -10318
$ws.Range("H126").Value = 1200
$ws.Range("I126").Value = 1200
$ws.Range("K126").Value = 3600
$ws.Range("M126").Value = -1130

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = ""
$ws.Range("N20").Value = ""
$ws.Range("H45").Value = 10000
$ws.Range("I45").Value = 10000
$ws.Range("K45").Value = 10000
$ws.Range("M45").Value = -9593
$ws.Range("H61").Value = 770.5
$ws.Range("I61").Value = 770.5
$ws.Range("K61").Value = 770.5
$ws.Range("M61").Value = -568.5
$ws.Range("H104").Value = 25327.8
$ws.Range("J104").Value = 25327.8
$ws.Range("L104").Value = 25327.8
$ws.Range("N104").Value = -32315.8
$ws.Range("H113").Value = 770.5
$ws.Range("I113").Value = 770.5
$ws.Range("K113").Value = 770.5
$ws.Range("M113").Value = 1399.5
$ws.Range("H132").Value = 750.75
$ws.Range("H135").Value = 55000
$ws.Range("J135").Value = 55000
$ws.Range("L135").Value = 55000
$ws.Range("N135").Value = -65140
$ws.Range("H136").Value = 21168.166
$ws.Range("I136").Value = 13000.75
$ws.Range("K136").Value = 39002.25
$ws.Range("M136").Value = -36452.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 275
$ws.Range("I132").Value = 275
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 825
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = 1705
$ws.Range("N132").Value = ""
$ws.Range("H136").Value = 2774.4666
$ws.Range("I136").Value = 1663.0769
$ws.Range("J136").Value = 9998.5
$ws.Range("K136").Value = 9998.5
$ws.Range("L136").Value = 29995.5
$ws.Range("M136").Value = -2439.2307
$ws.Range("N136").Value = -35095.5

Write-Output "applied Rafflesia_Profits market data refresh"
